# "Changes of File path in Address Book"
# Updates the ShipmentTrackNum (C) / PackageTrackNum (D) values for rows
# 2-22 on Sheet1 to a new batch of tracking numbers, and flips the Reject
# status in Q3 from FAIL to PASS.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [object]$Range,
        [string]$Value
    )
    # Force the cell to be written as text (shared string) rather than
    # being auto-coerced to a Number, without leaving a residual
    # NumberFormat/style applied to the cell (the source file keeps these
    # cells on the default "Normal" style).
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("C2") "320018612013"

# Row 3
Set-TextValue $ws.Range("C3") "320018612024"
Set-TextValue $ws.Range("Q3") "PASS"

# Row 4
Set-TextValue $ws.Range("C4") "320018612057"

# Row 5
Set-TextValue $ws.Range("C5") "320018612079"
Set-TextValue $ws.Range("D5") "320018612079"

# Row 6
Set-TextValue $ws.Range("C6") "320018612116"
Set-TextValue $ws.Range("D6") "320018612116"

# Row 7
Set-TextValue $ws.Range("C7") "320018612138"
Set-TextValue $ws.Range("D7") "320018612138"

# Row 8
Set-TextValue $ws.Range("C8") "320018612160"

# Row 9
Set-TextValue $ws.Range("C9") "320018612182"

# Row 10
Set-TextValue $ws.Range("C10") "320018612219"

# Row 11
Set-TextValue $ws.Range("C11") "320018612230"

# Row 12
Set-TextValue $ws.Range("C12") "320018612274"

# Row 13
Set-TextValue $ws.Range("C13") "320018612296"
Set-TextValue $ws.Range("D13") "320018612296"

# Row 14
Set-TextValue $ws.Range("C14") "320018612322"
Set-TextValue $ws.Range("D14") "320018612322"

# Row 15
Set-TextValue $ws.Range("C15") "320018612344"
Set-TextValue $ws.Range("D15") "320018612344"

# Row 16
Set-TextValue $ws.Range("C16") "320018612377"
Set-TextValue $ws.Range("D16") "320018612377"

# Row 17
Set-TextValue $ws.Range("C17") "320018612399"
Set-TextValue $ws.Range("D17") "320018612399"

# Row 18
Set-TextValue $ws.Range("C18") "320018612436"

# Row 19
Set-TextValue $ws.Range("C19") "320018612458"

# Row 20
Set-TextValue $ws.Range("C20") "320018612480"

# Row 21
Set-TextValue $ws.Range("C21") "320018612506"

# Row 22
Set-TextValue $ws.Range("C22") "320018612539"
